# Apply corrected IFRS metric values (rows 2-9) per "error solve ifrs list" fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "AA2" = 108.47
    "AB2" = 1490.8
    "AC2" = -34
    "AD2" = -170.04
    "AE2" = 16459
    "AF2" = 0.35
    "AG2" = 100
    "AH2" = 1.75
    "AI2" = -273.95
    "AJ2" = 24646734
    "D2" = 10441
    "E2" = 104
    "F2" = 104
    "G2" = 7
    "H2" = -8
    "I2" = -8
    "K2" = 7766
    "L2" = 4041
    "M2" = 3725
    "N2" = 3725
    "P2" = 246
    "Q2" = 557
    "R2" = -425
    "S2" = -92
    "T2" = 152
    "U2" = 405
    "V2" = 2188
    "W2" = 1
    "X2" = -0.08
    "Y2" = -0.22
    "Z2" = -0.1
    "AA3" = 57.09
    "AB3" = 1626.48
    "AC3" = 1458
    "AD3" = 5.72
    "AE3" = 17873
    "AF3" = 0.47
    "AG3" = 250
    "AH3" = 3
    "AI3" = 15.75
    "AJ3" = 24646734
    "D3" = 8897
    "E3" = 557
    "F3" = 557
    "G3" = 474
    "H3" = 359
    "I3" = 359
    "K3" = 6355
    "L3" = 2309
    "M3" = 4045
    "N3" = 4045
    "P3" = 246
    "Q3" = 1118
    "R3" = 747
    "S3" = -1427
    "T3" = 52
    "U3" = 1066
    "V3" = 784
    "W3" = 6.26
    "X3" = 4.04
    "Y3" = 9.25
    "Z3" = 5.09
    "AA4" = 52.1
    "AB4" = 1765.39
    "AC4" = 1622
    "AD4" = 5.86
    "AE4" = 19488
    "AF4" = 0.49
    "AG4" = 330
    "AH4" = 3.47
    "AI4" = 18.57
    "AJ4" = 24646734
    "D4" = 8979
    "E4" = 563
    "F4" = 563
    "G4" = 524
    "H4" = 398
    "I4" = 400
    "J4" = -1
    "K4" = 6696
    "L4" = 2294
    "M4" = 4402
    "N4" = 4386
    "O4" = 16
    "P4" = 246
    "Q4" = 409
    "R4" = -129
    "S4" = -208
    "T4" = 50
    "U4" = 359
    "V4" = 638
    "W4" = 6.27
    "X4" = 4.44
    "Y4" = 9.49
    "Z4" = 6.11
    "AA5" = 63.6
    "AB5" = 1859.12
    "AC5" = 1251
    "AD5" = 7.57
    "AE5" = 20446
    "AF5" = 0.46
    "AG5" = 300
    "AH5" = 3.17
    "AI5" = 21.9
    "AJ5" = 24646734
    "D5" = 12285
    "E5" = 471
    "F5" = 471
    "G5" = 400
    "H5" = 296
    "I5" = 308
    "J5" = -12
    "K5" = 7536
    "L5" = 2930
    "M5" = 4607
    "N5" = 4602
    "O5" = 5
    "P5" = 246
    "Q5" = 933
    "R5" = -440
    "S5" = -153
    "T5" = 101
    "U5" = 832
    "V5" = 581
    "W5" = 3.84
    "X5" = 2.41
    "Y5" = 6.86
    "Z5" = 4.16
    "AA6" = 51.09
    "AB6" = 1767.59
    "AC6" = -695
    "AD6" = -8.039999999999999
    "AE6" = 19356
    "AF6" = 0.29
    "AG6" = 150
    "AH6" = 2.68
    "AI6" = -19.69
    "AJ6" = 24646734
    "D6" = 11355
    "E6" = 25
    "F6" = 25
    "G6" = -180
    "H6" = -178
    "I6" = -171
    "K6" = 6571
    "L6" = 2222
    "M6" = 4349
    "N6" = 4351
    "P6" = 246
    "Q6" = -166
    "R6" = -138
    "S6" = -109
    "T6" = 173
    "U6" = -339
    "V6" = 545
    "W6" = 0.22
    "X6" = -1.57
    "Y6" = -3.83
    "Z6" = -2.52
    "AA7" = 47.29
    "AC7" = 1258
    "AD7" = 4.61
    "AE7" = 20551
    "AF7" = 0.28
    "AG7" = 150
    "AH7" = 2.59
    "AI7" = 11.93
    "D7" = 10730
    "E7" = 430
    "G7" = 390
    "H7" = 310
    "I7" = 310
    "K7" = 6790
    "L7" = 2180
    "M7" = 4610
    "N7" = 4620
    "P7" = 250
    "Q7" = 560
    "R7" = -340
    "S7" = -230
    "T7" = 170
    "W7" = 4.01
    "X7" = 2.89
    "Y7" = 6.91
    "Z7" = 4.64
    "AA8" = 44.08
    "AC8" = 1298
    "AD8" = 4.23
    "AE8" = 21841
    "AF8" = 0.25
    "AG8" = 150
    "AH8" = 2.73
    "AI8" = 11.55
    "D8" = 10670
    "E8" = 460
    "G8" = 430
    "H8" = 320
    "I8" = 320
    "K8" = 7060
    "L8" = 2160
    "M8" = 4900
    "N8" = 4910
    "P8" = 250
    "Q8" = 530
    "R8" = -530
    "S8" = -60
    "T8" = 170
    "W8" = 4.31
    "X8" = 3
    "Y8" = 6.72
    "Z8" = 4.62
    "AA9" = 42.31
    "AC9" = 1339
    "AD9" = 4.1
    "AE9" = 23131
    "AF9" = 0.24
    "AG9" = 150
    "AH9" = 2.73
    "AI9" = 11.2
    "D9" = 10710
    "E9" = 460
    "G9" = 440
    "H9" = 330
    "I9" = 330
    "K9" = 7400
    "L9" = 2200
    "M9" = 5200
    "N9" = 5200
    "P9" = 250
    "Q9" = 450
    "R9" = -390
    "S9" = -60
    "T9" = 170
    "W9" = 4.29
    "X9" = 3.08
    "Y9" = 6.53
    "Z9" = 4.56
}

foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}

# Cells removed entirely in the target (no longer populated)
$cellsToClear = @("J2", "O2", "J3", "O3", "U7", "U8", "U9")
foreach ($cellRef in $cellsToClear) {
    $ws.Range($cellRef).ClearContents()
}

